# Daily cryptocurrency price/volume refresh (GitHub Actions scheduled update).
# Updates Price (D) and Volume(1h) (E) columns for each coin row, and re-ranks
# rows 49-51 (Stellar / InjectiveProtocol / Bittensor swapped new rank order).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.616.64"
$ws.Range("E2").Value = "  -3.27%  "
# Row 3
$ws.Range("D3").Value = "2.968.01"
$ws.Range("E3").Value = "  -4.82%  "
# Row 4
$ws.Range("E4").Value = "  -0.15%  "
# Row 5
$ws.Range("D5").Value = "'543.25"
$ws.Range("E5").Value = "  -3.65%  "
# Row 6
$ws.Range("D6").Value = "'151.71"
$ws.Range("E6").Value = "  -4.81%  "
# Row 7
$ws.Range("E7").Value = "  -0.09%  "
# Row 8
$ws.Range("D8").Value = "'0.573"
$ws.Range("E8").Value = "  +2.24%  "
# Row 9
$ws.Range("D9").Value = "2.977.61"
$ws.Range("E9").Value = "  -4.79%  "
# Row 10
$ws.Range("D10").Value = "'0.113"
$ws.Range("E10").Value = "  -1.16%  "
# Row 11
$ws.Range("D11").Value = "'6.13"
$ws.Range("E11").Value = "  -6.04%  "
# Row 12
$ws.Range("E12").Value = "  -1.58%  "
# Row 13
$ws.Range("D13").Value = "3.486.46"
$ws.Range("E13").Value = "  -5.26%  "
# Row 14
$ws.Range("E14").Value = "  -2.59%  "
# Row 15
$ws.Range("D15").Value = "61.684.78"
$ws.Range("E15").Value = "  -3.43%  "
# Row 16
$ws.Range("E16").Value = "  -3.68%  "
# Row 17
$ws.Range("D17").Value = "2.974.60"
$ws.Range("E17").Value = "  -5.20%  "
# Row 18
$ws.Range("E18").Value = "  -3.47%  "
# Row 19
$ws.Range("E19").Value = "  +0.18%  "
# Row 20
$ws.Range("D20").Value = "'382.76"
$ws.Range("E20").Value = "  -3.52%  "
# Row 21
$ws.Range("D21").Value = "'12.02"
$ws.Range("E21").Value = "  -3.32%  "
# Row 22
$ws.Range("E22").Value = "  -4.79%  "
# Row 23
$ws.Range("E23").Value = "  +0.18%  "
# Row 24
$ws.Range("D24").Value = "'65.74"
$ws.Range("E24").Value = "  -2.98%  "
# Row 25
$ws.Range("D25").Value = "'0.470"
$ws.Range("E25").Value = "  -1.62%  "
# Row 26
$ws.Range("D26").Value = "3.092.17"
$ws.Range("E26").Value = "  -5.96%  "
# Row 27
$ws.Range("E27").Value = "  -2.64%  "
# Row 28
$ws.Range("E28").Value = "  -0.36%  "
# Row 29
$ws.Range("D29").Value = "0.0₃0940"
$ws.Range("E29").Value = "  -5.35%  "
# Row 30
$ws.Range("D30").Value = "'8.30"
$ws.Range("E30").Value = "  -3.94%  "
# Row 31
$ws.Range("E31").Value = "  +0.01%  "
# Row 32
$ws.Range("D32").Value = "'1.72"
$ws.Range("E32").Value = "  -3.64%  "
# Row 33
$ws.Range("D33").Value = "'20.49"
$ws.Range("E33").Value = "  -2.03%  "
# Row 34
$ws.Range("D34").Value = "'160.87"
$ws.Range("E34").Value = "  +2.77%  "
# Row 35
$ws.Range("D35").Value = "'4.64"
$ws.Range("E35").Value = "  -2.02%  "
# Row 36
$ws.Range("D36").Value = "'5.97"
$ws.Range("E36").Value = "  -3.29%  "
# Row 37
$ws.Range("E37").Value = "  -2.11%  "
# Row 38
$ws.Range("E38").Value = "  -3.45%  "
# Row 39
$ws.Range("E39").Value = "  -4.89%  "
# Row 40
$ws.Range("D40").Value = "'3.92"
$ws.Range("E40").Value = "  -2.24%  "
# Row 41
$ws.Range("D41").Value = "2.410.50"
$ws.Range("E41").Value = "  -8.66%  "
# Row 42
$ws.Range("D42").Value = "'37.27"
$ws.Range("E42").Value = "  -2.47%  "
# Row 43
$ws.Range("D43").Value = "'22.20"
$ws.Range("E43").Value = "  -4.87%  "
# Row 44
$ws.Range("E44").Value = "  -2.71%  "
# Row 45
$ws.Range("D45").Value = "'0.0595"
$ws.Range("E45").Value = "  -1.64%  "
# Row 46
$ws.Range("D46").Value = "'0.0249"
$ws.Range("E46").Value = "  -1.03%  "
# Row 47
$ws.Range("E47").Value = "  -0.02%  "
# Row 48
$ws.Range("D48").Value = "'5.03"
$ws.Range("E48").Value = "  -7.16%  "
# Row 49
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'19.87"
$ws.Range("E49").Value = "  -4.14%  "
# Row 50
$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").Value = "'269.43"
$ws.Range("E50").Value = "  -4.75%  "
# Row 51
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.0955"
$ws.Range("E51").Value = "  -1.23%  "

# The quote-prefix trick above can tag a cell with a text-number-format
# style; reset those cells back to the workbook default style so no
# extraneous formatting is introduced.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
